# Git & GitHub.pptx - reorder edit
#
# The author moved the slide titled "History of Version Control System"
# (originally slide #3) down so that it now sits right after
# "Basic flow to create git repository" (i.e. it becomes slide #10, just
# before "File status lifecycle"). All of the slides that used to sit
# between them (Git Features, Prerequisites, How to add repository,
# Git three stage architecture, Basic command, Basic flow to create git
# repository) shift up by one position; none of their own content changes.

$p = $ppt.ActivePresentation

$s = $p.Slides.Item(3)

# sanity check - make sure we are moving the right slide
# ("History of Version Control System")
$title = $s.Shapes.Item(1).TextFrame.TextRange.Text

$s.MoveTo(10)
